$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-17 Monday" "2025-11-18 Tuesday"

Replace-Text "95×82=7790" "57×39=2223"
Replace-Text "47×49=2303" "36×48=1728"
Replace-Text "67×28=1876" "13×65=845"
Replace-Text "25×79=1975" "87×53=4611"
Replace-Text "26×80=2080" "90×17=1530"

Replace-Text "50×88=4400" "81×89=7209"
Replace-Text "89×45=4005" "42×69=2898"
Replace-Text "52×54=2808" "71×97=6887"
Replace-Text "35×81=2835" "37×23=851"
Replace-Text "68×65=4420" "93×84=7812"

Replace-Text "22×81=1782" "45×68=3060"
Replace-Text "18×75=1350" "20×64=1280"
Replace-Text "39×51=1989" "35×22=770"
Replace-Text "27×54=1458" "88×66=5808"
Replace-Text "27×87=2349" "96×32=3072"

Replace-Text "74×69=5106" "15×79=1185"
Replace-Text "72×58=4176" "30×79=2370"
Replace-Text "32×50=1600" "58×90=5220"
Replace-Text "53×25=1325" "54×44=2376"
Replace-Text "41×47=1927" "55×57=3135"

Replace-Text "29×46=1334" "25×79=1975"
Replace-Text "22×89=1958" "43×12=516"
Replace-Text "42×83=3486" "39×16=624"
Replace-Text "55×87=4785" "28×86=2408"
Replace-Text "63×88=5544" "43×74=3182"
